$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 330, shifting rows 330:395 down to 331:396.
$ws.Rows.Item(330).Insert()

# Populate the new row 330 with a new weekly price record for
# "Femacal de La Calera" / "Apio" (same static attributes as the record
# that used to occupy row 330, but with a new date and price figures).
$ws.Range("A330").Value = 3
$ws.Range("B330").Value = "Femacal de La Calera"
$ws.Range("C330").Value = "Coquimbo"
$ws.Range("D330").Value = 44711
$ws.Range("E330").Value = 5
$ws.Range("F330").Value = 100112017
$ws.Range("G330").Value = "Apio"
$ws.Range("H330").Value = "Americana (o)"
$ws.Range("I330").Value = "Primera"
$ws.Range("J330").Value = 230
$ws.Range("K330").Value = 9000
$ws.Range("L330").Value = 9500
$ws.Range("M330").Value = 9239
$ws.Range("N330").Value = "$/docena de matas"
$ws.Range("O330").Value = "Pan de Az" + [char]0x00FA + "car"
$ws.Range("P330").Value = 1540
$ws.Range("Q330").Value = 6
$ws.Range("R330").Value = "Hortaliza"
